# feat: add pacjent, recepcjonistka, wizyta, wykonane badania controllers
#
# Fills in the "Piotr Bistyga" K:M column block (Data / Plik / Linie) for
# rows 18-21 with four newly-tracked files, and moves the active selection
# to K22 the way the author's last save left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18-21 all share the same "date" style already used by column N
# (style index 3, numFmtId 14 / m/d/yyyy). Copy formats from the sibling
# N column cell on each row so we reuse the existing style instead of
# Excel minting a brand-new cellXfs entry for "m/d/yyyy".
function Set-DateCell($row) {
    $srcFmt = $ws.Cells.Item($row, 14)   # N column, already styled as a date
    $dst = $ws.Cells.Item($row, 11)      # K column
    $dst.Value = 45793
    $srcFmt.Copy()
    $dst.PasteSpecial(-4122)             # xlPasteFormats
}

Set-DateCell 18
Set-DateCell 19
Set-DateCell 20
Set-DateCell 21

# L column: file name; M column: line count
$ws.Cells.Item(18, 12).Value = "WykonaneBadaniaController.cs"
$ws.Cells.Item(18, 13).Value = 35

$ws.Cells.Item(19, 12).Value = "PacjentController.cs"
$ws.Cells.Item(19, 13).Value = 43

$ws.Cells.Item(20, 12).Value = "RecepcjonistkaController.cs"
$ws.Cells.Item(20, 13).Value = 35

$ws.Cells.Item(21, 12).Value = "WizytaController.cs"
$ws.Cells.Item(21, 13).Value = 39

# Selection moves to K22, matching the author's last save.
[void]$ws.Range("K22").Select()

$excel.CutCopyMode = 0
